$wb = $excel.ActiveWorkbook

# "Activity Log - Part 2" is the 2nd worksheet in the workbook.
$ws = $wb.Worksheets.Item(2)

# Row 55
$ws.Range("D55").Value = 0.2951388888888889
$ws.Range("E55").Value = 0.3125

# Row 56
$ws.Range("D56").Value = 0.3125
$ws.Range("E56").Value = 0.31944444444444448

# Row 57
$ws.Range("D57").Value = 0.31944444444444448
$ws.Range("E57").Value = 0.3298611111111111

# Row 58
$ws.Range("D58").Value = 0.3298611111111111
$ws.Range("E58").Value = 0.34722222222222227

# Row 59 (also gains the last-4-digits / date values that were previously blank)
$ws.Range("B59").Value = 6977
$ws.Range("C59").Value = 43937
$ws.Range("D59").Value = 0.34722222222222227
$ws.Range("E59").Value = 0.3611111111111111

# Row 60 (also gains the last-4-digits / date values that were previously blank; no note text)
$ws.Range("B60").Value = 6977
$ws.Range("C60").Value = 43937
$ws.Range("D60").Value = 0.3611111111111111
$ws.Range("E60").Value = 0.36805555555555558

# The note text (column G) is entered in the same order the author typed it
# in, which is why the new shared-string table entries land in that order
# rather than row order (55, 56, 58, 59, 57).
$ws.Range("G55").Value = "Captured raw diagrams of functional waveforms from ModelSim for LogicUnit.vhd. - DONE"
$ws.Range("G56").Value = "Captured VHDL interface diagrams of all entities except ArithUnit.vhd (waiting on feedback from team members). -DONE"
$ws.Range("G58").Value = "Did not like an entity naming convenction so renamed it. Recompiled and recaptured VHDL interface, RTL and Post-Fit Diagrams. -DONE"
$ws.Range("G59").Value = "Captured raw diagrams of timing waveforms from ModelSim for LogicUnit.vhd. Exported transcript and summary files to Documentation folder for LogicUnit.vhd. -DONE"
$ws.Range("G57").Value = "Captured RTL and Post-Fit diagrams of LogicUnit and LogicGates. -DONE"

# Move the saved selection on this sheet to match the author's final cursor
# position (the workbook was left scrolled to show row 44 with B71 selected).
$ws.Activate()
$ws.Range("B71").Select()
